$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44658
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Especial"
$ws.Cells.Item(2, 13).Value = 280
$ws.Cells.Item(2, 14).Value = 21600
$ws.Cells.Item(2, 15).Value = 21600
$ws.Cells.Item(2, 16).Value = 21600
$ws.Cells.Item(2, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(2, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 19).Value = 1200
$ws.Cells.Item(2, 20).Value = 18
# Row 3
$ws.Cells.Item(3, 4).Value = 44658
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 330
$ws.Cells.Item(3, 14).Value = 16200
$ws.Cells.Item(3, 15).Value = 16200
$ws.Cells.Item(3, 16).Value = 16200
$ws.Cells.Item(3, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 900
$ws.Cells.Item(3, 20).Value = 18
# Row 4
$ws.Cells.Item(4, 4).Value = 44658
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 220
$ws.Cells.Item(4, 14).Value = 14400
$ws.Cells.Item(4, 15).Value = 14400
$ws.Cells.Item(4, 16).Value = 14400
$ws.Cells.Item(4, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 800
$ws.Cells.Item(4, 20).Value = 18
# Row 5
$ws.Cells.Item(5, 4).Value = 44309
$ws.Cells.Item(5, 11).Value = "Wonderfull"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 40
$ws.Cells.Item(5, 14).Value = 18000
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 18000
$ws.Cells.Item(5, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(5, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 19).Value = 1200
$ws.Cells.Item(5, 20).Value = 15
# Row 6
$ws.Cells.Item(6, 4).Value = 44309
$ws.Cells.Item(6, 11).Value = "Wonderfull"
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 70
$ws.Cells.Item(6, 14).Value = 15000
$ws.Cells.Item(6, 15).Value = 15000
$ws.Cells.Item(6, 16).Value = 15000
$ws.Cells.Item(6, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(6, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(6, 19).Value = 1000
$ws.Cells.Item(6, 20).Value = 15
# Row 7
$ws.Cells.Item(7, 4).Value = 44664
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Especial"
$ws.Cells.Item(7, 13).Value = 300
$ws.Cells.Item(7, 14).Value = 21600
$ws.Cells.Item(7, 15).Value = 21600
$ws.Cells.Item(7, 16).Value = 21600
$ws.Cells.Item(7, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(7, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 19).Value = 1200
$ws.Cells.Item(7, 20).Value = 18
# Row 8
$ws.Cells.Item(8, 4).Value = 44664
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 250
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 18
# Row 9
$ws.Cells.Item(9, 4).Value = 44664
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 250
$ws.Cells.Item(9, 14).Value = 16000
$ws.Cells.Item(9, 15).Value = 16000
$ws.Cells.Item(9, 16).Value = 16000
$ws.Cells.Item(9, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(9, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 19).Value = 889
$ws.Cells.Item(9, 20).Value = 18
# Row 10
$ws.Cells.Item(10, 4).Value = 44685
$ws.Cells.Item(10, 11).Value = "Wonderfull"
$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 350
$ws.Cells.Item(10, 14).Value = 21000
$ws.Cells.Item(10, 15).Value = 21000
$ws.Cells.Item(10, 16).Value = 21000
$ws.Cells.Item(10, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 1167
$ws.Cells.Item(10, 20).Value = 18
# Row 11
$ws.Cells.Item(11, 4).Value = 44685
$ws.Cells.Item(11, 11).Value = "Wonderfull"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 330
$ws.Cells.Item(11, 14).Value = 15000
$ws.Cells.Item(11, 15).Value = 15000
$ws.Cells.Item(11, 16).Value = 15000
$ws.Cells.Item(11, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 833
$ws.Cells.Item(11, 20).Value = 18
# Row 12
$ws.Cells.Item(12, 4).Value = 44685
$ws.Cells.Item(12, 11).Value = "Wonderfull"
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 280
$ws.Cells.Item(12, 14).Value = 10000
$ws.Cells.Item(12, 15).Value = 10000
$ws.Cells.Item(12, 16).Value = 10000
$ws.Cells.Item(12, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 556
$ws.Cells.Item(12, 20).Value = 18
# Row 13
$ws.Cells.Item(13, 4).Value = 44687
$ws.Cells.Item(13, 11).Value = "Wonderfull"
$ws.Cells.Item(13, 12).Value = "Especial"
$ws.Cells.Item(13, 13).Value = 220
$ws.Cells.Item(13, 14).Value = 21000
$ws.Cells.Item(13, 15).Value = 21000
$ws.Cells.Item(13, 16).Value = 21000
$ws.Cells.Item(13, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(13, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(13, 19).Value = 1167
$ws.Cells.Item(13, 20).Value = 18
# Row 14
$ws.Cells.Item(14, 4).Value = 44687
$ws.Cells.Item(14, 11).Value = "Wonderfull"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 250
$ws.Cells.Item(14, 14).Value = 15000
$ws.Cells.Item(14, 15).Value = 15000
$ws.Cells.Item(14, 16).Value = 15000
$ws.Cells.Item(14, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(14, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(14, 19).Value = 833
$ws.Cells.Item(14, 20).Value = 18
# Row 15
$ws.Cells.Item(15, 4).Value = 44687
$ws.Cells.Item(15, 11).Value = "Wonderfull"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 280
$ws.Cells.Item(15, 14).Value = 10000
$ws.Cells.Item(15, 15).Value = 10000
$ws.Cells.Item(15, 16).Value = 10000
$ws.Cells.Item(15, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(15, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(15, 19).Value = 556
$ws.Cells.Item(15, 20).Value = 18
# Row 16
$ws.Cells.Item(16, 4).Value = 44649
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 220
$ws.Cells.Item(16, 14).Value = 21600
$ws.Cells.Item(16, 15).Value = 21600
$ws.Cells.Item(16, 16).Value = 21600
$ws.Cells.Item(16, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(16, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 19).Value = 1200
$ws.Cells.Item(16, 20).Value = 18
# Row 17
$ws.Cells.Item(17, 4).Value = 44649
$ws.Cells.Item(17, 11).Value = "Sin especificar"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 250
$ws.Cells.Item(17, 14).Value = 16200
$ws.Cells.Item(17, 15).Value = 16200
$ws.Cells.Item(17, 16).Value = 16200
$ws.Cells.Item(17, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(17, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 19).Value = 900
$ws.Cells.Item(17, 20).Value = 18
# Row 18
$ws.Cells.Item(18, 4).Value = 44649
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 180
$ws.Cells.Item(18, 14).Value = 14400
$ws.Cells.Item(18, 15).Value = 14400
$ws.Cells.Item(18, 16).Value = 14400
$ws.Cells.Item(18, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(18, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 19).Value = 800
$ws.Cells.Item(18, 20).Value = 18
# Row 19
$ws.Cells.Item(19, 4).Value = 44285
$ws.Cells.Item(19, 11).Value = "Wonderfull"
$ws.Cells.Item(19, 12).Value = "Especial"
$ws.Cells.Item(19, 13).Value = 40
$ws.Cells.Item(19, 14).Value = 18000
$ws.Cells.Item(19, 15).Value = 18000
$ws.Cells.Item(19, 16).Value = 18000
$ws.Cells.Item(19, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(19, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 19).Value = 1200
$ws.Cells.Item(19, 20).Value = 15
# Row 20
$ws.Cells.Item(20, 4).Value = 44285
$ws.Cells.Item(20, 11).Value = "Wonderfull"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 90
$ws.Cells.Item(20, 14).Value = 15000
$ws.Cells.Item(20, 15).Value = 15000
$ws.Cells.Item(20, 16).Value = 15000
$ws.Cells.Item(20, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(20, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(20, 19).Value = 1000
$ws.Cells.Item(20, 20).Value = 15
# Row 21
$ws.Cells.Item(21, 4).Value = 44285
$ws.Cells.Item(21, 11).Value = "Wonderfull"
$ws.Cells.Item(21, 12).Value = "Segunda"
$ws.Cells.Item(21, 13).Value = 75
$ws.Cells.Item(21, 14).Value = 12000
$ws.Cells.Item(21, 15).Value = 12000
$ws.Cells.Item(21, 16).Value = 12000
$ws.Cells.Item(21, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(21, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(21, 19).Value = 800
$ws.Cells.Item(21, 20).Value = 15
# Row 22
$ws.Cells.Item(22, 4).Value = 44644
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Especial"
$ws.Cells.Item(22, 13).Value = 180
$ws.Cells.Item(22, 14).Value = 18000
$ws.Cells.Item(22, 15).Value = 18000
$ws.Cells.Item(22, 16).Value = 18000
$ws.Cells.Item(22, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(22, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 19).Value = 1200
$ws.Cells.Item(22, 20).Value = 15
# Row 23
$ws.Cells.Item(23, 4).Value = 44644
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 220
$ws.Cells.Item(23, 14).Value = 13500
$ws.Cells.Item(23, 15).Value = 13500
$ws.Cells.Item(23, 16).Value = 13500
$ws.Cells.Item(23, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(23, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 19).Value = 900
$ws.Cells.Item(23, 20).Value = 15
# Row 24
$ws.Cells.Item(24, 4).Value = 44644
$ws.Cells.Item(24, 11).Value = "Sin especificar"
$ws.Cells.Item(24, 12).Value = "Segunda"
$ws.Cells.Item(24, 13).Value = 290
$ws.Cells.Item(24, 14).Value = 12000
$ws.Cells.Item(24, 15).Value = 12000
$ws.Cells.Item(24, 16).Value = 12000
$ws.Cells.Item(24, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(24, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 19).Value = 800
$ws.Cells.Item(24, 20).Value = 15
# Row 25
$ws.Cells.Item(25, 4).Value = 44678
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 290
$ws.Cells.Item(25, 14).Value = 15000
$ws.Cells.Item(25, 15).Value = 15000
$ws.Cells.Item(25, 16).Value = 15000
$ws.Cells.Item(25, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(25, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(25, 19).Value = 1000
$ws.Cells.Item(25, 20).Value = 15
# Row 26
$ws.Cells.Item(26, 4).Value = 44678
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 220
$ws.Cells.Item(26, 14).Value = 12000
$ws.Cells.Item(26, 15).Value = 12000
$ws.Cells.Item(26, 16).Value = 12000
$ws.Cells.Item(26, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(26, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(26, 19).Value = 800
$ws.Cells.Item(26, 20).Value = 15
# Row 27
$ws.Cells.Item(27, 4).Value = 44305
$ws.Cells.Item(27, 11).Value = "Wonderfull"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 50
$ws.Cells.Item(27, 14).Value = 18000
$ws.Cells.Item(27, 15).Value = 18000
$ws.Cells.Item(27, 16).Value = 18000
$ws.Cells.Item(27, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(27, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(27, 19).Value = 1200
$ws.Cells.Item(27, 20).Value = 15
# Row 28
$ws.Cells.Item(28, 4).Value = 44305
$ws.Cells.Item(28, 11).Value = "Wonderfull"
$ws.Cells.Item(28, 12).Value = "Segunda"
$ws.Cells.Item(28, 13).Value = 60
$ws.Cells.Item(28, 14).Value = 15000
$ws.Cells.Item(28, 15).Value = 15000
$ws.Cells.Item(28, 16).Value = 15000
$ws.Cells.Item(28, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(28, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(28, 19).Value = 1000
$ws.Cells.Item(28, 20).Value = 15
